# Rename the display names of the inline logo pictures embedded in the
# document's headers/footers:
#   - footer1.xml Pearson logo:  image1.png -> image2.png
#   - footer2.xml Pearson logo:  image1.png -> image2.png
#   - header2.xml BTec logo:     image2.jpg -> image1.jpg
#
# These pictures live inside header/footer stories, so they are reached
# through Sections(1).Headers()/Footers() rather than the main
# Document.Range()/InlineShapes collection. We fetch each InlineShape via
# the paragraph that actually contains the <w:drawing> (the last paragraph
# of each header/footer story) so the object handle resolves cleanly.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-LogoInlineShape($story, $newName) {
    if ($story.Exists) {
        $range = $story.Range
        $count = $range.Paragraphs.Count
        $para = $range.Paragraphs($count)
        $shapes = $para.Range.InlineShapes
        if ($shapes.Count -ge 1) {
            $shape = $shapes.Item(1)
            $shape.Name = $newName
        }
    }
}

# BTEC logo (header2.xml) : image2.jpg -> image1.jpg
Rename-LogoInlineShape $sec.Headers(2) "image1.jpg"

# Pearson logo (footer1.xml) : image1.png -> image2.png
Rename-LogoInlineShape $sec.Footers(1) "image2.png"

# Pearson logo (footer2.xml) : image1.png -> image2.png
Rename-LogoInlineShape $sec.Footers(2) "image2.png"
